$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering date range) ---
$ws.Range("A8").Value = "Volume 31   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/8/2024  Through  4/14/2024"

# --- Crime-complaint data table updates (rows 14-30) ---
# F14: becomes the text placeholder "0" (style s="14")
$donor = $ws.Range("E14")
$donor.Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Formula = "=""0"""
$ws.Range("F14").Copy()
$ws.Range("F14").PasteSpecial(-4163)
$ws.Range("M14").Value = -85.714285714285
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -42.857142857142
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -84.615384615384
$ws.Range("C16").Value = 4
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -64
$ws.Range("I16").Value = 36
$ws.Range("J16").Value = 65
$ws.Range("K16").Value = -44.615384615384
$ws.Range("L16").Value = -7.692307692307
$ws.Range("M16").Value = -51.351351351351
$ws.Range("N16").Value = -89.565217391304
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -38.461538461538
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -36.666666666666
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 81
$ws.Range("K17").Value = -1.234567901234
$ws.Range("L17").Value = 1.265822784810
$ws.Range("M17").Value = -2.439024390243
$ws.Range("N17").Value = -62.441314553990
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 31
$ws.Range("J18").Value = 53
$ws.Range("K18").Value = -41.509433962264
$ws.Range("L18").Value = -50.793650793650
$ws.Range("M18").Value = -38
$ws.Range("N18").Value = -83.597883597883
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 13
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -43.478260869565
$ws.Range("I19").Value = 68
$ws.Range("J19").Value = 96
$ws.Range("K19").Value = -29.166666666666
$ws.Range("L19").Value = -33.980582524271
$ws.Range("M19").Value = -11.688311688311
$ws.Range("N19").Value = -29.166666666666
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 125
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = -3.703703703703
$ws.Range("L20").Value = -46.938775510204
$ws.Range("M20").Value = -3.703703703703
$ws.Range("N20").Value = -86.170212765957
$ws.Range("C21").Value = 23
$ws.Range("E21").Value = -17.857142857142
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = -34.736842105263
$ws.Range("I21").Value = 246
$ws.Range("J21").Value = 330
$ws.Range("K21").Value = -25.454545454545
$ws.Range("L21").Value = -28.488372093023
$ws.Range("M21").Value = -24.307692307692
$ws.Range("N21").Value = -76.901408450704
$ws.Range("D22").Value = 3
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -42.857142857142
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -55.555555555555
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 27
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = 3.846153846153
$ws.Range("L23").Value = 58.823529411764
$ws.Range("M23").Value = 12.5
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -35.294117647058
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = -9.459459459459
$ws.Range("I24").Value = 213
$ws.Range("J24").Value = 236
$ws.Range("K24").Value = -9.745762711864
$ws.Range("L24").Value = 8.121827411167
$ws.Range("M24").Value = 19.662921348314
$ws.Range("D25").Value = 4
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -66.666666666666
$ws.Range("J25").Value = 86
$ws.Range("K25").Value = -59.302325581395
$ws.Range("L25").Value = -20.454545454545
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -50
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -9.523809523809
$ws.Range("I26").Value = 101
$ws.Range("J26").Value = 122
$ws.Range("K26").Value = -17.213114754098
$ws.Range("L26").Value = -5.607476635514
$ws.Range("M26").Value = -46.842105263157
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 12
$ws.Range("K27").Value = -41.666666666666
$ws.Range("L27").Value = -36.363636363636
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = 0
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("I28").Value = 12
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = 140
$ws.Range("L28").Value = 33.333333333333
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 4
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = -75
$ws.Range("N29").Value = -92.727272727272
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = '#,##0'
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = -40
$ws.Range("L30").Value = -50
$ws.Range("M30").Value = -76.923076923076
$ws.Range("N30").Value = -94.117647058823

$excel.CutCopyMode = $false
